# Applies the "Add files via upload" commit to the colour-reference sheet:
#   - D264 gets recoloured from white (#ffffff) to red (#d20a11)
#   - four new rows (368-371) are appended for BECCS / CCS / DAC / Biochar
#     "negative emissions" colour entries (EN / FR / EN label + colour code)
#
# Column D carries a colour hex string whose cell fill is shaded to match,
# so for every colour cell we copy both value *and* format from an existing
# donor cell that already uses the exact colour we need - this reuses the
# existing style record instead of minting a near-duplicate one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D264: "#ffffff" (white, style 1) -> "#d20a11" (red, style 8) ---------
$ws.Range("D263").Copy($ws.Range("D264"))

# --- Row 368: BECCS (negative emissions) ----------------------------------
$ws.Range("A368").Value = "BECCS (negative)"
$ws.Range("B368").Value = "BECSC (émissions négatives)"
$ws.Range("C368").Value = "BECCS (negative emissions)"
$ws.Range("D363").Copy($ws.Range("D368"))

# --- Row 369: CCS / CSC ----------------------------------------------------
$ws.Range("A369").Value = "CCS"
$ws.Range("B369").Value = "CSC"
$ws.Range("C369").Value = "CCS"
$ws.Range("D17").Copy($ws.Range("D369"))

# --- Row 370: DAC (negative emissions) -------------------------------------
$ws.Range("A370").Value = "DAC (negative)"
$ws.Range("B370").Value = "EDA (émissions négatives)"
$ws.Range("C370").Value = "DAC (negative emissions)"
$ws.Range("D83").Copy($ws.Range("D370"))

# --- Row 371: Biochar (negative emissions) ----------------------------------
$ws.Range("A371").Value = "Biochar (negative)"
$ws.Range("B371").Value = "Biocharbon (émissions négatives)"
$ws.Range("C371").Value = "Biochar (negative emissions)"
$ws.Range("D275").Copy($ws.Range("D371"))

# Leave the cursor on D264, matching where the author ended up after editing
# the colour of that cell.
[void]$ws.Range("D264").Select()
